$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the range that was edited (matches resulting selection in the file)
$ws.Range("B1:E3").Select()

# Row 1 header values updated
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data updated; C2 and E2 cleared (removed)
$ws.Range("B2").Value = 24.099583937430452
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 15.025160356548326
$ws.Range("E2").ClearContents()

# Row 3 data updated; C3 cleared (removed), D3 newly populated
$ws.Range("B3").Value = 22.555286873196565
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 15.166599691751969
$ws.Range("E3").Value = 30.483408351428295
